# core/text.xlsx -- "fullfil some client functions and examples"
#
# 1. sheet1 ("第n批次德育分申请表（收集结果）"): duplicate the one data
#    row (row 2) twice, into rows 3 and 4, so the collected-results sheet
#    has three sample submissions instead of one.
# 2. sheet2 ("mapping"): add a "column name" (B) next to the existing
#    Python-identifier (A) / description (C) columns, add a new row 15
#    with the "submitter (auto)" field, and add a formula column (D)
#    that renders each A/B pair as a Python dict literal fragment
#    ("{'key':'label'},") for easy copy/paste into client code.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# Sheet 1: duplicate the sample data row (row 2) into rows 3 and 4.
# ---------------------------------------------------------------------------
$ws1.Range("A2:M2").Copy()
$ws1.Range("A3:M3").PasteSpecial()
$ws1.Range("A2:M2").Copy()
$ws1.Range("A4:M4").PasteSpecial()
$excel.CutCopyMode = $false

$ws1.Rows.Item(3).RowHeight = 39
$ws1.Rows.Item(4).RowHeight = 39

# ---------------------------------------------------------------------------
# Sheet 2: new "column name" values in B3:B15.
# ---------------------------------------------------------------------------
$ws2.Range("B3").Value = "学号（必填）"
$ws2.Range("B4").Value = "姓名（必填）"
$ws2.Range("B5").Value = "申请类型（必填）"
$ws2.Range("B6").Value = "申请分数（必填）"
$ws2.Range("B7").Value = "材料附件（必填）"
$ws2.Range("B8").Value = "活动描述（必填）"
$ws2.Range("B9").Value = "活动时间（必填）"
$ws2.Range("B10").Value = "备注"
$ws2.Range("B11").Value = "审核状态"
$ws2.Range("B12").Value = "审核员"
$ws2.Range("B13").Value = "审核备注"
$ws2.Range("B14").Value = "审核日期"
$ws2.Range("B15").Value = "提交者（自动）"

# Column D: dict-literal fragment built from A (key) and B (label).
$ws2.Range("D3").Formula = "=""{'""&A3&""'""&"":""&""'""&B3&""'""&""}""&"","""
$ws2.Range("D4:D14").Formula = "=""{'""&A4&""'""&"":""&""'""&B4&""'""&""}""&"","""

$ws1.Activate()
$excel.ActiveWindow.FreezePanes = $false
$ws1.Range("D1").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws1.Range("E15").Select()

$ws2.Activate()
$ws2.Range("D3:D14").Select()

$wb.Save()
